$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# like "1.000" or "0.7751" are preserved exactly as typed, matching the
# original inlineStr cell type/content instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.968.61'
$ws.Range("E2").Value = '  +0.53%  '

$ws.Range("D3").Value = '1.892.85'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '0.7751'
$ws.Range("E5").Value = '  -0.57%  '

$ws.Range("D6").Value = '243.92'
$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '0.3131'
$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").Value = '25.82'
$ws.Range("E9").Value = '  +2.41%  '

$ws.Range("D10").Value = '0.07260'
$ws.Range("E10").Value = '  +1.70%  '

$ws.Range("D11").Value = '0.08671'
$ws.Range("E11").Value = '  +7.38%  '

$ws.Range("D12").Value = '1.986.86'
$ws.Range("E12").Value = '  +5.19%  '

$ws.Range("D13").Value = '0.7732'
$ws.Range("E13").Value = '  +1.68%  '

$ws.Range("D14").Value = '5.420'
$ws.Range("E14").Value = '  -0.53%  '

$ws.Range("D15").Value = '94.51'
$ws.Range("E15").Value = '  +2.74%  '

$ws.Range("D16").Value = '6.220'
$ws.Range("E16").Value = '  +1.29%  '

$ws.Range("D17").Value = '30.171.52'
$ws.Range("E17").Value = '  +1.25%  '

$ws.Range("E18").Value = '  +0.32%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '245.71'
$ws.Range("E19").Value = '  +1.09%  '

$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.298.37'
$ws.Range("E20").Value = '  +8.44%  '

$ws.Range("E21").Value = '  +1.82%  '

$ws.Range("D22").Value = '8.171'
$ws.Range("E22").Value = '  +0.90%  '

$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("D25").Value = '0.1599'
$ws.Range("E25").Value = '  -1.70%  '

$ws.Range("D26").Value = '9.547'
$ws.Range("E26").Value = '  +1.72%  '

$ws.Range("D27").Value = '162.87'
$ws.Range("E27").Value = '  +0.38%  '

$ws.Range("D28").Value = '18.86'
$ws.Range("E28").Value = '  +0.99%  '

$ws.Range("D29").Value = '2.049'

$ws.Range("E30").Value = '  +1.60%  '

$ws.Range("D31").Value = '1.547'

$ws.Range("D32").Value = '4.533'
$ws.Range("E32").Value = '  +1.58%  '

$ws.Range("D33").Value = '4.133'
$ws.Range("E33").Value = '  +0.95%  '

$ws.Range("E34").Value = '  -1.24%  '

$ws.Range("D35").Value = '1.251'
$ws.Range("E35").Value = '  -0.81%  '

$ws.Range("D36").Value = '0.7556'
$ws.Range("E36").Value = '  +1.92%  '

$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  +0.60%  '

$ws.Range("E38").Value = '  +2.55%  '

$ws.Range("D39").Value = '0.01980'
$ws.Range("E39").Value = '  +3.50%  '

$ws.Range("D40").Value = '2.785'
$ws.Range("E40").Value = '  +0.29%  '

$ws.Range("D41").Value = '0.4526'
$ws.Range("E41").Value = '  +2.81%  '

$ws.Range("D42").Value = '73.80'
$ws.Range("E42").Value = '  +0.44%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '6.063'
$ws.Range("E43").Value = '  +3.93%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.093.10'
$ws.Range("E44").Value = '  -4.06%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.8543'
$ws.Range("E45").Value = '  +0.56%  '

$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.214.20'
$ws.Range("E46").Value = '  +8.44%  '

$ws.Range("D47").Value = '0.9997'
$ws.Range("E47").Value = '  -0.10%  '

$ws.Range("D48").Value = '103.44'
$ws.Range("E48").Value = '  +0.01%  '

$ws.Range("E49").Value = '  +1.29%  '

$ws.Range("D50").Value = '7.631'
$ws.Range("E50").Value = '  +2.78%  '

$ws.Range("D51").Value = '9.864'
$ws.Range("E51").Value = '  -0.33%  '

# Restore the default (Normal) cell style on column D so no stray number-format
# style is left applied to the cells (matches original workbook formatting).
$ws.Range("D2:D51").Style = "Normal"
